$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "_old" -> "_FV2310", "_new" -> "_FV2404"
#    (row 1, columns A:J are the "_old" set, L:U are the "_new" set, K="diff")
# ---------------------------------------------------------------------------
$headersFV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (pane split below row 1)
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the data range into a native Excel Table ("Table1") so the
#    sheet gets a tableParts/<table> definition with the same 21 columns
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

Write-Output "done"
